{"js": "// Locate the (single) changelog table in the document body.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// --- Change 1: the row that currently reads \"01/03/2022\" becomes \"03/01/2021\" ---\n// (row index 4, column index 0 - zero based: header is row 0)\nconst dateCellOld = table.getCell(4, 0);\nconst dateParaOld = dateCellOld.body.paragraphs.getFirst();\ndateParaOld.insertText(\"03/01/2021\", Word.InsertLocation.replace);\n\n// --- Change 2: fill in the next (previously empty) changelog row ---\n// row index 5, columns: 0=Ng\u00e0y l\u1eadp, 1=M\u00f4 t\u1ea3 thay \u0111\u1ed5i, 2=Phi\u00ean b\u1ea3n, 3=Ng\u01b0\u1eddi l\u1eadp\nconst newRowDate = table.getCell(5, 0);\nconst newRowDesc = table.getCell(5, 1);\nconst newRowVersion = table.getCell(5, 2);\nconst newRowAuthor = table.getCell(5, 3);\n\nnewRowDate.body.paragraphs.getFirst().insertText(\"16/05/2021\", Word.InsertLocation.replace);\nnewRowDesc.body.paragraphs.getFirst().insertText(\n  \"S\u1eeda l\u1ed7i \u0111i\u1ec3m kh\u00f4ng t\u0103ng khi \u0111\u1ea1t \u0111\u1ebfn 1600 \u0111i\u1ec3m\",\n  Word.InsertLocation.replace\n);\nnewRowVersion.body.paragraphs.getFirst().insertText(\"0.2.2\", Word.InsertLocation.replace);\nnewRowAuthor.body.paragraphs.getFirst().insertText(\"Nguy\u1ec5n Th\u00e0nh Long\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# --- Change 1: the row that currently reads \"01/03/2022\" becomes \"03/01/2021\" ---\n# (COM tables/cells are 1-based; header is row 1, so this is row 5, column 1)\n$table.Cell(5, 1).Range.Text = \"03/01/2021\"\n\n# --- Change 2: fill in the next (previously empty) changelog row ---\n# row 6: 1=Ng\u00e0y l\u1eadp, 2=M\u00f4 t\u1ea3 thay \u0111\u1ed5i, 3=Phi\u00ean b\u1ea3n, 4=Ng\u01b0\u1eddi l\u1eadp\n$table.Cell(6, 1).Range.Text = \"16/05/2021\"\n$table.Cell(6, 2).Range.Text = \"S\u1eeda l\u1ed7i \u0111i\u1ec3m kh\u00f4ng t\u0103ng khi \u0111\u1ea1t \u0111\u1ebfn 1600 \u0111i\u1ec3m\"\n$table.Cell(6, 3).Range.Text = \"0.2.2\"\n$table.Cell(6, 4).Range.Text = \"Nguy\u1ec5n Th\u00e0nh Long\"\n"}
